$p = $ppt.ActivePresentation
Write-Host "Designs.Count before:" $p.Designs.Count
try {
  $d = $p.Designs.Add()
  Write-Host "Designs.Add() ->" $d
} catch {
  Write-Host "Designs.Add failed:" $_.Exception.Message
}
Write-Host "Designs.Count after Add:" $p.Designs.Count
